# Split melatonin worksheets: MLT; PLB; ALL
# Update the "No. Hours" entry for the week of row 15 (B15:C15 = 42912-42918)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E15").Value = 22
